$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.312.34'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '2.587.78'
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.44'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.78'
$ws.Range('E6').Value = '  +3.16%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.542'
$ws.Range('E9').Value = '  +2.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.89'
$ws.Range('E10').Value = '  -0.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0816'
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.54'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '2.985.43'
$ws.Range('E13').Value = '  +2.38%  '
$ws.Range('E14').Value = '  -1.95%  '
$ws.Range('D15').Value = '2.602.43'
$ws.Range('E15').Value = '  +3.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.29'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.848'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '43.342.07'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.87'
$ws.Range('E19').Value = '  +3.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.56'
$ws.Range('E20').Value = '  -2.76%  '
$ws.Range('D21').Value = '0.0₃0965'
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.68'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '255.31'
$ws.Range('E23').Value = '  +1.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.99'
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.09'
$ws.Range('E25').Value = '  +3.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.35'
$ws.Range('E26').Value = '  +2.04%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('E28').Value = '  +1.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.22'
$ws.Range('E29').Value = '  +0.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.36'
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.88'
$ws.Range('E31').Value = '  -2.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.36'
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.44'
$ws.Range('E33').Value = '  +4.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.17'
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0810'
$ws.Range('E35').Value = '  +2.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.70'
$ws.Range('E36').Value = '  +3.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.71'
$ws.Range('E37').Value = '  -2.02%  '
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('E39').Value = '  +8.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.119'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.65'
$ws.Range('E41').Value = '  -3.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.98'
$ws.Range('E42').Value = '  +5.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0305'
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.27'
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('D46').Value = '2.018.44'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.00'
$ws.Range('E47').Value = '  +2.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '83.61'
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.836.08'
$ws.Range('E49').Value = '  +2.37%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.51'
$ws.Range('E50').Value = '  +2.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.194'
$ws.Range('E51').Value = '  +2.64%  '
